$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. Rename sheet "INPUT_PATH" -> "INPUT_SETTING"
# ----------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("INPUT_PATH")
$wsInput.Name = "INPUT_SETTING"

# ----------------------------------------------------------------------
# 2. Update INPUT_SETTING cell contents
#    (DV_AZUL path now points at the FILTER_UL sheet, output file names
#     renamed, and a new Output UL row is populated)
# ----------------------------------------------------------------------
$wsInput.Range("B7").Value = "refer to FILTER_UL sheet"
$wsInput.Range("B8").Value = "test_filetrad"
$wsInput.Range("B9").Value = "test_fileul"
$wsInput.Range("A6").Select()

# ----------------------------------------------------------------------
# 3. FILTER_TRAD - only the active selection moved
# ----------------------------------------------------------------------
$wsTrad = $wb.Worksheets.Item("FILTER_TRAD")
$wsTrad.Range("G6").Select()

# ----------------------------------------------------------------------
# 4. FILTER_UL - add path_dv / path_uvsg columns, and two data rows
# ----------------------------------------------------------------------
$wsUL = $wb.Worksheets.Item("FILTER_UL")

# Insert two new columns: one right after "run_name" (for path_dv) and
# one right after "path_rafm" (for path_uvsg)
$wsUL.Columns.Item(2).Insert()
$wsUL.Columns.Item(4).Insert()

$wsUL.Range("B1").Value = "path_dv"
$wsUL.Range("D1").Value = "path_uvsg"

$wsUL.Range("A2").Value = "run4"
$wsUL.Range("B2").Value = "D:\1. IRCS Automation\Control 3 DEV\Source\DV_AZUL_Stat_Con_2Q25.xlsx"
$wsUL.Range("C2").Value = "D:\1. IRCS Automation\Control 3 DEV\Source\Data_Extraction_run4UL_Con.xlsx"
$wsUL.Range("D2").Value = "D:\1. IRCS Automation\Control 3 DEV\Source\Data_Extraction_run4UVSG.xlsx"
$wsUL.Range("E2").Value = 16233
$wsUL.Range("J2").Value = "SH,UL,PI"

$wsUL.Range("A3").Value = "run5"
$wsUL.Range("B3").Value = "D:\1. IRCS Automation\Control 3 DEV\Source\DV_AZUL_Stat_Con_2Q25.xlsx"
$wsUL.Range("C3").Value = "D:\1. IRCS Automation\Control 3 DEV\Source\Data_Extraction_run4UL_Con.xlsx"
$wsUL.Range("D3").Value = "D:\1. IRCS Automation\Control 3 DEV\Source\Data_Extraction_run4UVSG.xlsx"
$wsUL.Range("E3").Value = 16234
$wsUL.Range("J3").Value = "SH,UL,PI"

# ----------------------------------------------------------------------
# 5. VARIABLE_DEF - only the active selection moved
# ----------------------------------------------------------------------
$wsVar = $wb.Worksheets.Item("VARIABLE_DEF")
$wsVar.Range("C9").Select()

# ----------------------------------------------------------------------
# 6. FILTER_UL becomes the active / selected sheet last, so it "sticks"
#    as the workbook's active tab.
# ----------------------------------------------------------------------
$wsUL.Activate()
$wsUL.Range("C3").Select()
